$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.115.90'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.251.75'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.06'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '185.06'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.47%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.600'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.66%  '
$ws.Range('E9').Value = '  -2.73%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.62'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.418'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.66%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.812.54'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.137'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E14').Value = '  -2.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '68.112.98'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.34%  '
$ws.Range('E16').Value = '  -0.41%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.227.14'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.75%  '
$ws.Range('E18').Value = '  -0.27%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.51'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '395.46'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.33%  '
$ws.Range('E21').Value = '  -0.58%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.50'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.31%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.517'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.79%  '
$ws.Range('E25').Value = '  -0.55%  '
$ws.Range('E26').Value = '  +3.17%  '
$ws.Range('E27').Value = '  -3.16%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.64'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.72%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.79'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.17%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.06'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.27'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.14%  '
$ws.Range('E35').Value = '  -4.92%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '162.16'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.23%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.92'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.62%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.817'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.63'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.60%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '26.47'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.53'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.06%  '
$ws.Range('E42').Value = '  -4.52%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '41.03'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.21%  '
$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0686'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '25.19'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.71%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.609.07'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.93%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '338.55'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0280'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.16%  '
$ws.Range('E49').Value = '  +3.24%  '
$ws.Range('E50').Value = '  -0.69%  '
$ws.Range('B51').Value = 'Arweave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '31.21'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.61%  '
